$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("L2").Value = 1.29
$ws.Range("AF2").Value = 25
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 46
$ws.Range("AO2").Value = 21

# Row 3
$ws.Range("S3").Value = 3.15

# Row 4
$ws.Range("Q4").Value = 1.68
$ws.Range("R4").Value = 1.46
$ws.Range("S4").Value = 2.5
$ws.Range("AD4").Value = 16.5

# Row 5
$ws.Range("AB5").Value = 990
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 990

# Row 7
$ws.Range("U7").Value = 1.7
$ws.Range("AN7").Value = 15
